$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.188.66'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.858.03'
$ws.Range('E3').Value = '  -1.23%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.02'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7023'
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3111'
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07772'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.15'
$ws.Range('E10').Value = '  -4.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07977'
$ws.Range('E11').Value = '  -4.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.858.21'
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '93.36'
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.166'
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6952'
$ws.Range('E15').Value = '  -3.30%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.344'
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.160.42'
$ws.Range('E17').Value = '  -0.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008269'
$ws.Range('E18').Value = '  -4.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '250.45'
$ws.Range('E19').Value = '  +3.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.110.87'
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.06'
$ws.Range('E21').Value = '  -1.30%  '
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.495'
$ws.Range('E23').Value = '  -4.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1546'
$ws.Range('E25').Value = '  -2.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.964'
$ws.Range('E26').Value = '  -1.06%  '
$ws.Range('E27').Value = '  -2.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.78'
$ws.Range('E28').Value = '  +1.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.495'
$ws.Range('E29').Value = '  -0.97%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.290'
$ws.Range('E30').Value = '  -2.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.251'
$ws.Range('E31').Value = '  -2.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.206'
$ws.Range('E32').Value = '  +0.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05246'
$ws.Range('E33').Value = '  -2.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.876'
$ws.Range('E34').Value = '  -3.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7424'
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.153'
$ws.Range('E36').Value = '  -2.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.710'
$ws.Range('E37').Value = '  +0.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01859'
$ws.Range('E38').Value = '  -1.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.243.12'
$ws.Range('E39').Value = '  -3.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.732'
$ws.Range('E40').Value = '  -0.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.208'
$ws.Range('E41').Value = '  -6.00%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8952'
$ws.Range('E42').Value = '  -2.51%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '110.43'
$ws.Range('E43').Value = '  -1.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '70.82'
$ws.Range('E44').Value = '  -5.05%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('E46').Value = '  +0.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.009.23'
$ws.Range('E47').Value = '  -1.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5181'
$ws.Range('E48').Value = '  -0.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.776'
$ws.Range('E49').Value = '  -1.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.404'
$ws.Range('E50').Value = '  -1.12%  '
$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4287'
$ws.Range('E51').Value = '  -2.21%  '
